$d = $word.ActiveDocument

$startPara = 0
$endPara = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
  $t = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]10)
  if ($t -eq "Icon und splashscreen" -and $startPara -eq 0) {
    $startPara = $i
  }
  if ($t -eq "/resources/splash.png" -and $startPara -ne 0 -and $i -gt $startPara) {
    $endPara = $i
  }
}

if ($startPara -eq 0 -or $endPara -eq 0) {
  throw ("Anchor paragraphs not found: start=" + $startPara + " end=" + $endPara)
}

$rangeStart = $d.Paragraphs($startPara).Range.Start
$rangeEnd = $d.Paragraphs($endPara).Range.End
$full = $d.Range($rangeStart, $rangeEnd)

$xmlFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="berschrift2"/></w:pPr><w:r><w:t>Abfrage per CLI</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>curl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -D - -X "POST" http://geoweb.zamg.ac.at/quakeapi/v02/getapikey -H "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Authorization</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: Basic cXVha2VhcGk6I3FrcCZtbGRuZyM=" -H "Content-Type: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>application</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>json</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>charset</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=utf-8"</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>{"apikey":"8415dcd4-e88e-11e6-a0a2-525401d06</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>b63"}</w:t></w:r></w:p><w:p><w:r><w:t>8415dcd4-e88e-11e6-a0a2-525401d06b63</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>curl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -D - -X "POST" -H "X-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>QuakeAPIKey</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: 8415dcd4-e88e-11e6-a0a2-525401d06b63" --</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>data-binary</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> @./test2.json "http://geoweb.zamg.ac.at/quakeapi/v02/message" -H "</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Authorization</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: Basic cXVha2VhcGk6I3FrcCZtbGRuZyM=" -H "Content-Type: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>application</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>json</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>charset</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=utf-8"</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="berschrift2"/></w:pPr><w:r><w:t xml:space="preserve">Icon und </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>splashscreen</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>resources</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/icon.png</w:t></w:r></w:p><w:p><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>resources</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/splash.png</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$full.InsertXML($xmlFrag)

Write-Output ("paragraphs now: " + $d.Paragraphs.Count)
